$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total "Valor Mora" summary figure (row 11)
$ws.Range("E11").Value = 444000

# Remove 4 rows from the detail table (16 rows -> 12 rows); deleting rows
# 27:30 (rather than 28:31) preserves the special "last row" styling that
# lived on the physical row 31 so it lands back on the new last row (27).
$ws.Rows("27:30").Delete()

# Rewrite the detail rows (now 16-27) with the refreshed data: two new
# workers (Miguel Angel Marsiglia Muñoz / Carmen Ortiz Angulo) covering
# periods 2206-2211, and Tatiana Martinez Vasquez now only for 2208-2211.
$ws.Range("C16").Value = "73151013"
$ws.Range("D16").Value = "MIGUEL ANGEL MARSIGLIA MUÑOZ"
$ws.Range("E16").Value = "2206"
$ws.Range("F16").Value = 40000

$ws.Range("C17").Value = "1051885383"
$ws.Range("D17").Value = "CARMEN ORTIZ ANGULO"
$ws.Range("E17").Value = "2206"
$ws.Range("F17").Value = 40000

$ws.Range("C18").Value = "73151013"
$ws.Range("D18").Value = "MIGUEL ANGEL MARSIGLIA MUÑOZ"
$ws.Range("E18").Value = "2207"
$ws.Range("F18").Value = 40000

$ws.Range("C19").Value = "1051885383"
$ws.Range("D19").Value = "CARMEN ORTIZ ANGULO"
$ws.Range("E19").Value = "2207"
$ws.Range("F19").Value = 40000

$ws.Range("C20").Value = "64743394"
$ws.Range("D20").Value = "TATIANA DEL SOCORRO MARTINEZ VASQUEZ"
$ws.Range("E20").Value = "2208"
$ws.Range("F20").Value = 28000

$ws.Range("C21").Value = "73151013"
$ws.Range("D21").Value = "MIGUEL ANGEL MARSIGLIA MUÑOZ"
$ws.Range("E21").Value = "2208"
$ws.Range("F21").Value = 40000

$ws.Range("C22").Value = "64743394"
$ws.Range("D22").Value = "TATIANA DEL SOCORRO MARTINEZ VASQUEZ"
$ws.Range("E22").Value = "2209"
$ws.Range("F22").Value = 40000

$ws.Range("C23").Value = "73151013"
$ws.Range("D23").Value = "MIGUEL ANGEL MARSIGLIA MUÑOZ"
$ws.Range("E23").Value = "2209"
$ws.Range("F23").Value = 40000

$ws.Range("C24").Value = "64743394"
$ws.Range("D24").Value = "TATIANA DEL SOCORRO MARTINEZ VASQUEZ"
$ws.Range("E24").Value = "2210"
$ws.Range("F24").Value = 40000

$ws.Range("C25").Value = "73151013"
$ws.Range("D25").Value = "MIGUEL ANGEL MARSIGLIA MUÑOZ"
$ws.Range("E25").Value = "2210"
$ws.Range("F25").Value = 40000

$ws.Range("C26").Value = "64743394"
$ws.Range("D26").Value = "TATIANA DEL SOCORRO MARTINEZ VASQUEZ"
$ws.Range("E26").Value = "2211"
$ws.Range("F26").Value = 28000

$ws.Range("C27").Value = "73151013"
$ws.Range("D27").Value = "MIGUEL ANGEL MARSIGLIA MUÑOZ"
$ws.Range("E27").Value = "2211"
$ws.Range("F27").Value = 28000
